$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.615.99"
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("D3").Value = "1.665.67"
$ws.Range("E3").Value = "  -4.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'215.64"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'0.509"
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'24.08"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.263"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("D12").Value = "1.902.64"
$ws.Range("E12").Value = "  -3.98%  "
$ws.Range("D13").Value = "1.659.42"
$ws.Range("E13").Value = "  -4.39%  "
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").Value = "'0.568"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'66.45"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "27.607.73"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "'241.88"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -3.70%  "
$ws.Range("D20").Value = "'7.67"
$ws.Range("E20").Value = "  -4.14%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("E23").Value = "  -4.11%  "
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").Value = "'147.05"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("D27").Value = "'16.38"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").Value = "'0.0504"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("D33").Value = "1.465.54"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.37"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.929"
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").Value = "'0.577"
$ws.Range("E38").Value = "  -5.07%  "
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'69.96"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("E41").Value = "  -4.19%  "
$ws.Range("E43").Value = "  -3.56%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.796"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.40"
$ws.Range("E45").Value = "  -5.64%  "
$ws.Range("D46").Value = "1.809.65"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "'1.74"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "'88.92"
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'7.92"
$ws.Range("E51").Value = "  -3.69%  "

# Some prices are plain numeric-looking strings (e.g. "215.64"). Assigning
# them directly would make Excel auto-convert to a Double (losing the exact
# text and introducing float rounding). We force them to text with a leading
# apostrophe, then paste-special just the number format from a blank, default-
# styled cell so the cells keep the workbook-default "General" style instead of
# picking up a quote-prefixed style.
$fmtSource = $ws.Range("Z1")
$fmtSource.Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("D51").PasteSpecial(-4122)
$excel.CutCopyMode = $false
